$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 31, pushing the existing rows 31-70 down to 32-71.
$ws.Rows("31:31").Insert()

# Populate the newly inserted row 31 with the new weekly record.
$ws.Cells.Item(31, 1).Value = 8
$ws.Cells.Item(31, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(31, 3).Value = "Coquimbo"
$ws.Cells.Item(31, 4).Value = (Get-Date -Year 1899 -Month 12 -Day 30 -Hour 0 -Minute 0 -Second 0).AddDays(44413)
$ws.Cells.Item(31, 5).Value = 4
$ws.Cells.Item(31, 6).Value = 100112001
$ws.Cells.Item(31, 7).Value = "Berenjena"
$ws.Cells.Item(31, 8).Value = "Sin especificar"
$ws.Cells.Item(31, 9).Value = "Primera"
$ws.Cells.Item(31, 10).Value = 640
$ws.Cells.Item(31, 11).Value = 12000
$ws.Cells.Item(31, 12).Value = 13000
$ws.Cells.Item(31, 13).Value = 12500
$ws.Cells.Item(31, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(31, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(31, 16).Value = 208
$ws.Cells.Item(31, 17).Value = 60
$ws.Cells.Item(31, 18).Value = "Hortaliza"
